# Insert a new row above current row 89, shifting rows 89:96 down to 90:97.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(89).Insert()

# Populate the newly inserted row 89 with the new weekly price record.
$ws.Range("A89").Value = 1
$ws.Range("B89").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C89").Value = "Arica y Parinacota"
$ws.Range("D89").Value = 44783
$ws.Range("E89").Value = 15
$ws.Range("F89").Value = "Fruta"
$ws.Range("G89").Value = 100102
$ws.Range("H89").Value = "Cítricos"
$ws.Range("I89").Value = 100102005
$ws.Range("J89").Value = "Naranja"
$ws.Range("K89").Value = "Navel Late"
$ws.Range("L89").Value = "Tercera"
$ws.Range("M89").Value = 250
$ws.Range("N89").Value = 500
$ws.Range("O89").Value = 550
$ws.Range("P89").Value = 525
$ws.Range("Q89").Value = "`$/kilo (en caja de 20 kilos)"
$ws.Range("R89").Value = "Región de Coquimbo"
$ws.Range("S89").Value = 525
$ws.Range("T89").Value = 1
